$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    under the H1 title at the top of the document.
# ---------------------------------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $metaPara = $findRng.Paragraphs(1)
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Near the bottom of the document, right before the "Prompt: ..." image
#    prompt paragraph, add a new bold paragraph repeating the page title, and
#    turn the "Prompt: ..." paragraph itself into the (former) meta
#    description text, keeping its italic run formatting.
# ---------------------------------------------------------------------------
$promptRng = $d.Content
$foundPrompt = $promptRng.Find.Execute("Prompt:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundPrompt) {
    $promptPara = $promptRng.Paragraphs(1)
    $targetRng = $promptPara.Range

    $titleText = "Play Creepy Carnival Free - Review of NoLimit City's Spooky Slot"
    $descText = "Explore the eerie circus show of NoLimit City's Creepy Carnival. Play for free and enjoy two unique features - the Free Spin and Star Spin mode."

    $xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
    $fragment = "<w:p $xmlNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>$titleText</w:t></w:r></w:p>" +
                "<w:p $xmlNs><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$descText</w:t></w:r></w:p>"

    [void]$targetRng.InsertXML($fragment)
}
